$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.443260550498962
$ws.Range("B1").Value = 3.311005353927612
$ws.Range("C1").Value = 5.392126083374023
$ws.Range("D1").Value = 7.282094955444336
$ws.Range("E1").Value = 1.001720786094666
